$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at row 43, shifting existing rows 43.. down by one
# (this pushes the "Broadband" label in A168 down to A169, etc.).
$ws.Rows.Item(43).Insert()

# Populate the new row's September details / date columns with the
# newest transaction.
$ws.Range("R43").Value = "debit"
$ws.Range("S43").Value = "2024-09-19 14:35:16"
